# Scheduled-runner style refresh of cached marketboard/profit figures across
# the per-job sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR). Values below are
# plain numeric snapshots (no formulas) being overwritten in place; a couple
# of rows also gain/lose a trailing HQ-profit cell as the source data shifts.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 275
$ws.Range("I33").Value = 272.4
$ws.Range("K33").Value = 272.4
$ws.Range("M33").Value = -43.39999999999998
$ws.Range("H74").Value = 3812.375
$ws.Range("J74").Value = 4150
$ws.Range("L74").Value = 4150
$ws.Range("N74").Value = -6022
$ws.Range("H77").Value = 3812.375
$ws.Range("J77").Value = 4150
$ws.Range("L77").Value = 20750
$ws.Range("N77").Value = -30110
$ws.Range("H86").Value = 2100.5
$ws.Range("I86").Value = 1937.875
$ws.Range("J86").Value = 2317.3333
$ws.Range("K86").Value = 1937.875
$ws.Range("L86").Value = 2317.3333
$ws.Range("M86").Value = -814.875
$ws.Range("N86").Value = -4563.3333
$ws.Range("H89").Value = 2100.5
$ws.Range("I89").Value = 1937.875
$ws.Range("J89").Value = 2317.3333
$ws.Range("K89").Value = 9689.375
$ws.Range("L89").Value = 11586.6665
$ws.Range("M89").Value = -4073.375
$ws.Range("N89").Value = -22818.6665
$ws.Range("H92").Value = 266.4
$ws.Range("I92").Value = 189.25
$ws.Range("K92").Value = 189.25
$ws.Range("M92").Value = 1058.75
$ws.Range("H132").Value = 1520.8032
$ws.Range("I132").Value = 1094.64
$ws.Range("J132").Value = 3457.9092
$ws.Range("K132").Value = 3283.92
$ws.Range("L132").Value = 10373.7276
$ws.Range("M132").Value = -753.9200000000001
$ws.Range("N132").Value = -15433.7276
$ws.Range("H138").Value = 2200859.8
$ws.Range("I138").Value = 11112842
$ws.Range("J138").Value = 3384.6711
$ws.Range("K138").Value = 33338526
$ws.Range("L138").Value = 10154.0133
$ws.Range("M138").Value = -33333386
$ws.Range("N138").Value = -20434.0133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2350
$ws.Range("I61").Value = 2350
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2350
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2138
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 1356.1305
$ws.Range("I74").Value = 1247.5454
$ws.Range("J74").Value = 1455.6666
$ws.Range("K74").Value = 1247.5454
$ws.Range("L74").Value = 1455.6666
$ws.Range("M74").Value = -373.5454
$ws.Range("N74").Value = -3203.6666
$ws.Range("H77").Value = 1356.1305
$ws.Range("I77").Value = 1247.5454
$ws.Range("J77").Value = 1455.6666
$ws.Range("K77").Value = 6237.727
$ws.Range("L77").Value = 7278.333000000001
$ws.Range("M77").Value = -1869.727
$ws.Range("N77").Value = -16014.333
$ws.Range("H136").Value = 2350
$ws.Range("I136").Value = 2350
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7050
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4500
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1334.8636
$ws.Range("I94").Value = 1030.4667
$ws.Range("K94").Value = 1030.4667
$ws.Range("M94").Value = -579.4666999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1523.6471
$ws.Range("I122").Value = 1314.5
$ws.Range("K122").Value = 3943.5
$ws.Range("M122").Value = -1493.5
$ws.Range("H134").Value = 942.6316
$ws.Range("I134").Value = 968.58826
$ws.Range("J134").Value = 722
$ws.Range("K134").Value = 2905.76478
$ws.Range("L134").Value = 2166
$ws.Range("M134").Value = -370.76478
$ws.Range("N134").Value = -7236
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3047.4666
$ws.Range("I102").Value = 2882.9092
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 2882.9092
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -1260.9092
$ws.Range("N102").Value = -6744
$ws.Range("H113").Value = 1082.9333
$ws.Range("I113").Value = 577.2
$ws.Range("J113").Value = 2094.4
$ws.Range("K113").Value = 577.2
$ws.Range("L113").Value = 2094.4
$ws.Range("M113").Value = 1592.8
$ws.Range("N113").Value = -6434.4
$ws.Range("H122").Value = 2475.8438
$ws.Range("I122").Value = 2705.5
$ws.Range("J122").Value = 1786.875
$ws.Range("K122").Value = 8116.5
$ws.Range("L122").Value = 5360.625
$ws.Range("M122").Value = -5666.5
$ws.Range("N122").Value = -10260.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9932
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 9932
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 9932
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -10156
$ws.Range("H5").Value = 18900
$ws.Range("I5").Value = 26000
$ws.Range("J5").Value = 14166.667
$ws.Range("K5").Value = 26000
$ws.Range("L5").Value = 14166.667
$ws.Range("M5").Value = -25887
$ws.Range("N5").Value = -14392.667
$ws.Range("H17").Value = 7269.6665
$ws.Range("I17").Value = 800
$ws.Range("J17").Value = 10504.5
$ws.Range("K17").Value = 800
$ws.Range("L17").Value = 10504.5
$ws.Range("M17").Value = -630
$ws.Range("N17").Value = -10844.5
$ws.Range("H40").Value = 12444.333
$ws.Range("I40").Value = 16666.5
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 16666.5
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -16530.5
$ws.Range("N40").Value = -4272
$ws.Range("H93").Value = 1628.4286
$ws.Range("I93").Value = 1480
$ws.Range("K93").Value = 1480
$ws.Range("M93").Value = -232
$ws.Range("H96").Value = 29194
$ws.Range("J96").Value = 29194
$ws.Range("L96").Value = 29194
$ws.Range("N96").Value = -34686
$ws.Range("H122").Value = 10420469
$ws.Range("I122").Value = 13161592
$ws.Range("J122").Value = 4200
$ws.Range("K122").Value = 39484776
$ws.Range("L122").Value = 12600
$ws.Range("M122").Value = -39482326
$ws.Range("N122").Value = -17500
$ws.Range("H136").Value = 1779.2646
$ws.Range("I136").Value = 1806.2903
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 5418.8709
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -2868.8709
$ws.Range("N136").Value = -9600
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 101939.9
$ws.Range("I81").Value = 126449.875
$ws.Range("J81").Value = 3900
$ws.Range("K81").Value = 252899.75
$ws.Range("L81").Value = 7800
$ws.Range("M81").Value = -251838.75
$ws.Range("N81").Value = -9922
$ws.Range("H84").Value = 101939.9
$ws.Range("I84").Value = 126449.875
$ws.Range("J84").Value = 3900
$ws.Range("K84").Value = 1264498.75
$ws.Range("L84").Value = 39000
$ws.Range("M84").Value = -1259194.75
$ws.Range("N84").Value = -49608
$ws.Range("H122").Value = 20838506
$ws.Range("I122").Value = 50004040
$ws.Range("J122").Value = 5981.4287
$ws.Range("K122").Value = 150012120
$ws.Range("L122").Value = 17944.2861
$ws.Range("M122").Value = -150009670
$ws.Range("N122").Value = -22844.2861
